# Updates the cryptos list data (prices & volume %) on sheet1 to the
# latest scrape. Two pairs of rows had their coin order swapped:
#   - Uniswap <-> WrappedEther (rows 17/18)
#   - Bittensor <-> InjectiveProtocol (rows 39/40)
# All D/E-column cells hold text (not real numbers/percentages), so every
# assignment below forces a text NumberFormat first to stop Excel's COM
# layer from auto-coercing number-looking strings (e.g. "587.67") into
# floating point values, which would silently lose the exact original
# formatting/precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.447.24"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "3.561.69"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.67"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.22"
$ws.Range("E6").Value = "  +3.01%  "
$ws.Range("D7").Value = "3.559.28"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E8").Value = "  +1.72%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  +10.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.646"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.18"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000311"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "4.130.36"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "70.480.79"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.575.12"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.77"
$ws.Range("E18").Value = "  +4.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.96"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "580.19"
$ws.Range("E20").Value = "  +9.03%  "
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.73"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.63"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.87"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.41"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.07"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.20"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.88"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  +24.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.23"
$ws.Range("E36").Value = "  +5.48%  "
$ws.Range("D37").Value = "3.782.20"
$ws.Range("E37").Value = "  +12.98%  "
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.06"
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "523.93"
$ws.Range("E40").Value = "  -4.20%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "0.0₃0787"
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("E43").Value = "  +5.70%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0455"
$ws.Range("E45").Value = "  +4.56%  "
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.20"
$ws.Range("E49").Value = "  +2.24%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.43"
$ws.Range("E51").Value = "  +7.60%  "
